# Update to levee setback and levee removal codes
# Mostly work updating the results plotting

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # mf_wb_dict
$ws2 = $wb.Worksheets.Item(2)   # owhm_wb_dict

# --- New worksheet: flopy_to_owhm ------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $last)
$ws3.Name = "flopy_to_owhm"

# header row (written B before A so "flopy" becomes shared-string #71 and
# "owhm" becomes #72, matching the source order of new strings)
$ws3.Range("B2").Value = "flopy"
$ws3.Range("A2").Value = "owhm"

# --- Add two new rows to owhm_wb_dict (GHB_IN / GHB_OUT) -------------------
# Column A values are new shared strings, B/C reuse existing ones.
$ws2.Range("A12").Value = "GHB_IN"
$ws2.Range("B12").Value = "black"
$ws2.Range("C12").Value = "Subsurface Inflow"

$ws2.Range("A13").Value = "GHB_OUT"
$ws2.Range("B13").Value = "black"
$ws2.Range("C13").Value = "Subsurface Outflow"

# the comment header for the new sheet
$ws3.Range("A1").Value = "# dictionary matching flopy wb summary to flopy zonebudget naming (GW_IN is a manual addition)"

# --- Storage-change row added to both owhm_wb_dict and flopy_to_owhm -------
$ws2.Range("A14").Value = "dSTORAGE_sum"
$ws2.Range("B14").Value = "black"
$ws2.Range("C14").Value = "Cumulative Storage Change"

# --- Fill in the rest of the flopy_to_owhm mapping table --------------------
$ws3.Range("A3").Value = "RCH_IN"
$ws3.Range("B3").Value = "FROM_RECHARGE"

$ws3.Range("A4").Value = "SFR_IN"
$ws3.Range("B4").Value = "FROM_STREAM_LEAKAGE"

$ws3.Range("A5").Value = "LAK_IN"
$ws3.Range("B5").Value = "FROM_LAKE_SEEPAGE"

$ws3.Range("A6").Value = "GW_IN"
$ws3.Range("B6").Value = "FROM_HEAD_DEP_BOUNDS"

$ws3.Range("A7").Value = "WEL_OUT"
$ws3.Range("B7").Value = "TO_WELLS"

$ws3.Range("A8").Value = "ET_OUT"
$ws3.Range("B8").Value = "TO_ET"

$ws3.Range("A9").Value = "SFR_OUT"
$ws3.Range("B9").Value = "TO_STREAM_LEAKAGE"

$ws3.Range("A10").Value = "GW_OUT"
$ws3.Range("B10").Value = "TO_HEAD_DEP_BOUNDS"

$ws3.Range("A11").Value = "LAK_OUT"
$ws3.Range("B11").Value = "TO_LAKE_SEEPAGE"

$ws3.Range("A12").Value = "GHB_IN"
$ws3.Range("B12").Value = "FROM_HEAD_DEP_BOUNDS"

$ws3.Range("A13").Value = "GHB_OUT"
$ws3.Range("B13").Value = "TO_HEAD_DEP_BOUNDS"

$ws3.Range("A14").Value = "dSTORAGE_sum"
$ws3.Range("B14").Value = "dSTORAGE_sum"

# --- column widths on the new sheet -----------------------------------------
$ws3.Columns.Item(1).ColumnWidth = 14.385416666666666
$ws3.Columns.Item(2).ColumnWidth = 24.276041666666668

# --- selections -------------------------------------------------------------
$ws1.Range("A18").Select()
$ws2.Range("A14").Select()
$ws3.Range("B14").Select()
